$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63. This shifts the previous rows 63..163
# down to 64..164, while leaving row 62's existing content untouched in
# place (still the "Sutil De Gase" record) and leaving the newly inserted
# row 63 blank.
$ws.Rows.Item(63).Insert()

# The target layout wants the record that used to live in row 62 to now
# occupy row 63 (unchanged), and row 62 to hold a brand-new record. So
# first clone row 62's current values down into the blank row 63...
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
foreach ($c in $cols) {
    $ws.Range("$c`63").Value2 = $ws.Range("$c`62").Value2
}

# ...then overwrite row 62 in place with the new record's data. Only the
# fields that actually differ from the old row 62 are touched; the rest
# (Mercado/Region/Codreg/Tipo/Producto/Categoria/Calidad/Unidad/Origen/Kg)
# stay exactly as they already are.
$ws.Range("D62").Value = 44498
$ws.Range("K62").Value = "Tahití"
$ws.Range("M62").Value = 250
$ws.Range("N62").Value = 36000
$ws.Range("O62").Value = 37000
$ws.Range("P62").Value = 36500
$ws.Range("S62").Value = 1521
